# Reinstate the APL and ANL carrier columns (parser revert) and rename
# YES -> CMA, matching the "temporarily reverted parser due to bug" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) YES -> CMA rename (same A:B columns, values keep their row-category text)
$ws.Range("A1").Value = "CMA picks"
$ws.Range("B1").Value = "CMA drops"

# 2) Make room for two new carriers (APL, ANL) right after CMA by inserting
#    6 blank columns at D:I. Everything from D onward (COS, HDM, HLC, ... ZIM)
#    shifts right by 6 columns, which reproduces the target layout exactly.
$ws.Columns("D:I").Insert()

# 3) Fill in the APL (D:E) and ANL (G:H) header + body cells.
#    Each carrier occupies two columns: "<CODE> picks" / "<CODE> drops" in
#    row 1, and the row's own size-category label in rows 2-7 (mirroring the
#    pattern used by every other carrier column pair on this sheet).
$rowLabels = @{ 2 = "20'"; 3 = "40DR"; 4 = "40DH"; 5 = "45DH"; 6 = "Reefer"; 7 = "Special" }

$ws.Range("D1").Value = "APL picks"
$ws.Range("E1").Value = "APL drops"
$ws.Range("G1").Value = "ANL picks"
$ws.Range("H1").Value = "ANL drops"

foreach ($r in 2..7) {
    $label = $rowLabels[$r]
    $ws.Cells.Item($r, 4).Value = $label   # D
    $ws.Cells.Item($r, 5).Value = $label   # E
    $ws.Cells.Item($r, 7).Value = $label   # G
    $ws.Cells.Item($r, 8).Value = $label   # H
}

# 4) HLC (now P:Q) picked up two extra "drops" data points for 40DH/45DH.
$ws.Range("Q4").Value = "40DH"
$ws.Range("Q5").Value = "45DH"
